$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Ghi chu" (Note) column E
$ws.Range("E2").Value = "Ghi chú"
$ws.Range("E3").Value = "Nhóm trưởng"
$ws.Range("E4").Value = "Thành viên"
$ws.Range("E5").Value = "Thành viên"
$ws.Range("E6").Value = "Thành viên"
$ws.Range("E7").Value = "Thành viên"

# Extra notes below the table
$ws.Range("B9").Value = "Đường dẫn thùng chứa:http://quanlythuvien5n.googlecode.com/svn/trunk/"
$ws.Range("B10").Value = "Đường dẫn đến project: http://code.google.com/p/quanlythuvien5n/"

$ws.Columns.Item(5).ColumnWidth = 14.36

$ws.Range("B14").Select() | Out-Null
